$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Vtn"
$ws.Cells.Item(2, 3).Value = "Itgb6"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = [double]"3"
$ws.Cells.Item(2, 6).Value = [double]"1"
$ws.Cells.Item(2, 7).Value = [double]"2.767552"
$ws.Cells.Item(2, 8).Value = [double]"8.302655999999999"
$ws.Cells.Item(2, 9).Value = [double]"0.04706493447833917"
$ws.Cells.Item(2, 10).Value = [double]"0.04706493447833917"
$ws.Cells.Item(2, 11).Value = [double]"3"
$ws.Cells.Item(2, 12).Value = [double]"1"
$ws.Cells.Item(2, 13).Value = [double]"0.2577576666666667"
$ws.Cells.Item(2, 14).Value = [double]"0.7732730000000001"
$ws.Cells.Item(2, 15).Value = [double]"0.03524815007985697"
$ws.Cells.Item(2, 16).Value = [double]"0.03524815007985697"
$ws.Cells.Item(2, 17).Value = [double]"0.7133577458986667"
$ws.Cells.Item(2, 18).Value = [double]"6.420219713088"
$ws.Cells.Item(2, 19).Value = [double]"0.001658951873991134"
$ws.Cells.Item(2, 20).Value = [double]"0.001658951873991134"

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Vtn"
$ws.Cells.Item(3, 3).Value = "Itgb6"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = [double]"3"
$ws.Cells.Item(3, 6).Value = [double]"1"
$ws.Cells.Item(3, 7).Value = [double]"2.767552"
$ws.Cells.Item(3, 8).Value = [double]"8.302655999999999"
$ws.Cells.Item(3, 9).Value = [double]"0.04706493447833917"
$ws.Cells.Item(3, 10).Value = [double]"0.04706493447833917"
$ws.Cells.Item(3, 11).Value = [double]"3"
$ws.Cells.Item(3, 12).Value = [double]"1"
$ws.Cells.Item(3, 13).Value = [double]"2.242708666666667"
$ws.Cells.Item(3, 14).Value = [double]"6.728126"
$ws.Cells.Item(3, 15).Value = [double]"0.3066885757089511"
$ws.Cells.Item(3, 16).Value = [double]"0.3066885757089511"
$ws.Cells.Item(3, 17).Value = [double]"6.206812855850666"
$ws.Cells.Item(3, 18).Value = [double]"55.86131570265599"
$ws.Cells.Item(3, 19).Value = [double]"0.01443427772099695"
$ws.Cells.Item(3, 20).Value = [double]"0.01443427772099695"

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Vtn"
$ws.Cells.Item(4, 3).Value = "Itgb6"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = [double]"3"
$ws.Cells.Item(4, 6).Value = [double]"1"
$ws.Cells.Item(4, 7).Value = [double]"2.767552"
$ws.Cells.Item(4, 8).Value = [double]"8.302655999999999"
$ws.Cells.Item(4, 9).Value = [double]"0.04706493447833917"
$ws.Cells.Item(4, 10).Value = [double]"0.04706493447833917"
$ws.Cells.Item(4, 11).Value = [double]"3"
$ws.Cells.Item(4, 12).Value = [double]"1"
$ws.Cells.Item(4, 13).Value = [double]"4.789377333333333"
$ws.Cells.Item(4, 14).Value = [double]"14.368132"
$ws.Cells.Item(4, 15).Value = [double]"0.6549434327891901"
$ws.Cells.Item(4, 16).Value = [double]"0.6549434327891902"
$ws.Cells.Item(4, 17).Value = [double]"13.25485081762133"
$ws.Cells.Item(4, 18).Value = [double]"119.293657358592"
$ws.Cells.Item(4, 19).Value = [double]"0.03082486975124176"
$ws.Cells.Item(4, 20).Value = [double]"0.03082486975124177"

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Vtn"
$ws.Cells.Item(5, 3).Value = "Itgb6"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = [double]"3"
$ws.Cells.Item(5, 6).Value = [double]"1"
$ws.Cells.Item(5, 7).Value = [double]"2.767552"
$ws.Cells.Item(5, 8).Value = [double]"8.302655999999999"
$ws.Cells.Item(5, 9).Value = [double]"0.04706493447833917"
$ws.Cells.Item(5, 10).Value = [double]"0.04706493447833917"
$ws.Cells.Item(5, 11).Value = [double]"1"
$ws.Cells.Item(5, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(5, 13).Value = [double]"0.02281433333333334"
$ws.Cells.Item(5, 14).Value = [double]"0.068443"
$ws.Cells.Item(5, 15).Value = [double]"0.003119841422001868"
$ws.Cells.Item(5, 16).Value = [double]"0.003119841422001868"
$ws.Cells.Item(5, 17).Value = [double]"0.06313985384533334"
$ws.Cells.Item(5, 18).Value = [double]"0.568258684608"
$ws.Cells.Item(5, 19).Value = [double]"0.0001468351321093264"
$ws.Cells.Item(5, 20).Value = [double]"0.0001468351321093264"

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Vtn"
$ws.Cells.Item(6, 3).Value = "Itgb6"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = [double]"3"
$ws.Cells.Item(6, 6).Value = [double]"1"
$ws.Cells.Item(6, 7).Value = [double]"21.05317333333333"
$ws.Cells.Item(6, 8).Value = [double]"63.15952"
$ws.Cells.Item(6, 9).Value = [double]"0.3580298485789791"
$ws.Cells.Item(6, 10).Value = [double]"0.3580298485789791"
$ws.Cells.Item(6, 11).Value = [double]"3"
$ws.Cells.Item(6, 12).Value = [double]"1"
$ws.Cells.Item(6, 13).Value = [double]"0.2577576666666667"
$ws.Cells.Item(6, 14).Value = [double]"0.7732730000000001"
$ws.Cells.Item(6, 15).Value = [double]"0.03524815007985697"
$ws.Cells.Item(6, 16).Value = [double]"0.03524815007985697"
$ws.Cells.Item(6, 17).Value = [double]"5.42661683432889"
$ws.Cells.Item(6, 18).Value = [double]"48.83955150896001"
$ws.Cells.Item(6, 19).Value = [double]"0.01261988983578032"
$ws.Cells.Item(6, 20).Value = [double]"0.01261988983578032"

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Vtn"
$ws.Cells.Item(7, 3).Value = "Itgb6"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = [double]"3"
$ws.Cells.Item(7, 6).Value = [double]"1"
$ws.Cells.Item(7, 7).Value = [double]"21.05317333333333"
$ws.Cells.Item(7, 8).Value = [double]"63.15952"
$ws.Cells.Item(7, 9).Value = [double]"0.3580298485789791"
$ws.Cells.Item(7, 10).Value = [double]"0.3580298485789791"
$ws.Cells.Item(7, 11).Value = [double]"3"
$ws.Cells.Item(7, 12).Value = [double]"1"
$ws.Cells.Item(7, 13).Value = [double]"2.242708666666667"
$ws.Cells.Item(7, 14).Value = [double]"6.728126"
$ws.Cells.Item(7, 15).Value = [double]"0.3066885757089511"
$ws.Cells.Item(7, 16).Value = [double]"0.3066885757089511"
$ws.Cells.Item(7, 17).Value = [double]"47.21613429550222"
$ws.Cells.Item(7, 18).Value = [double]"424.94520865952"
$ws.Cells.Item(7, 19).Value = [double]"0.1098036643219785"
$ws.Cells.Item(7, 20).Value = [double]"0.1098036643219785"

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Vtn"
$ws.Cells.Item(8, 3).Value = "Itgb6"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = [double]"3"
$ws.Cells.Item(8, 6).Value = [double]"1"
$ws.Cells.Item(8, 7).Value = [double]"21.05317333333333"
$ws.Cells.Item(8, 8).Value = [double]"63.15952"
$ws.Cells.Item(8, 9).Value = [double]"0.3580298485789791"
$ws.Cells.Item(8, 10).Value = [double]"0.3580298485789791"
$ws.Cells.Item(8, 11).Value = [double]"3"
$ws.Cells.Item(8, 12).Value = [double]"1"
$ws.Cells.Item(8, 13).Value = [double]"4.789377333333333"
$ws.Cells.Item(8, 14).Value = [double]"14.368132"
$ws.Cells.Item(8, 15).Value = [double]"0.6549434327891901"
$ws.Cells.Item(8, 16).Value = [double]"0.6549434327891902"
$ws.Cells.Item(8, 17).Value = [double]"100.8315911574044"
$ws.Cells.Item(8, 18).Value = [double]"907.48432041664"
$ws.Cells.Item(8, 19).Value = [double]"0.2344892980693105"
$ws.Cells.Item(8, 20).Value = [double]"0.2344892980693106"

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Vtn"
$ws.Cells.Item(9, 3).Value = "Itgb6"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = [double]"3"
$ws.Cells.Item(9, 6).Value = [double]"1"
$ws.Cells.Item(9, 7).Value = [double]"21.05317333333333"
$ws.Cells.Item(9, 8).Value = [double]"63.15952"
$ws.Cells.Item(9, 9).Value = [double]"0.3580298485789791"
$ws.Cells.Item(9, 10).Value = [double]"0.3580298485789791"
$ws.Cells.Item(9, 11).Value = [double]"1"
$ws.Cells.Item(9, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(9, 13).Value = [double]"0.02281433333333334"
$ws.Cells.Item(9, 14).Value = [double]"0.068443"
$ws.Cells.Item(9, 15).Value = [double]"0.003119841422001868"
$ws.Cells.Item(9, 16).Value = [double]"0.003119841422001868"
$ws.Cells.Item(9, 17).Value = [double]"0.4803141141511112"
$ws.Cells.Item(9, 18).Value = [double]"4.322827027360001"
$ws.Cells.Item(9, 19).Value = [double]"0.001116996351909756"
$ws.Cells.Item(9, 20).Value = [double]"0.001116996351909756"

$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Vtn"
$ws.Cells.Item(10, 3).Value = "Itgb6"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = [double]"3"
$ws.Cells.Item(10, 6).Value = [double]"1"
$ws.Cells.Item(10, 7).Value = [double]"34.97741266666667"
$ws.Cells.Item(10, 8).Value = [double]"104.932238"
$ws.Cells.Item(10, 9).Value = [double]"0.5948251867999219"
$ws.Cells.Item(10, 10).Value = [double]"0.5948251867999219"
$ws.Cells.Item(10, 11).Value = [double]"3"
$ws.Cells.Item(10, 12).Value = [double]"1"
$ws.Cells.Item(10, 13).Value = [double]"0.2577576666666667"
$ws.Cells.Item(10, 14).Value = [double]"0.7732730000000001"
$ws.Cells.Item(10, 15).Value = [double]"0.03524815007985697"
$ws.Cells.Item(10, 16).Value = [double]"0.03524815007985697"
$ws.Cells.Item(10, 17).Value = [double]"9.015696274997115"
$ws.Cells.Item(10, 18).Value = [double]"81.14126647497402"
$ws.Cells.Item(10, 19).Value = [double]"0.02096648745560261"
$ws.Cells.Item(10, 20).Value = [double]"0.02096648745560261"

$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Vtn"
$ws.Cells.Item(11, 3).Value = "Itgb6"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = [double]"3"
$ws.Cells.Item(11, 6).Value = [double]"1"
$ws.Cells.Item(11, 7).Value = [double]"34.97741266666667"
$ws.Cells.Item(11, 8).Value = [double]"104.932238"
$ws.Cells.Item(11, 9).Value = [double]"0.5948251867999219"
$ws.Cells.Item(11, 10).Value = [double]"0.5948251867999219"
$ws.Cells.Item(11, 11).Value = [double]"3"
$ws.Cells.Item(11, 12).Value = [double]"1"
$ws.Cells.Item(11, 13).Value = [double]"2.242708666666667"
$ws.Cells.Item(11, 14).Value = [double]"6.728126"
$ws.Cells.Item(11, 15).Value = [double]"0.3066885757089511"
$ws.Cells.Item(11, 16).Value = [double]"0.3066885757089511"
$ws.Cells.Item(11, 17).Value = [double]"78.44414652510979"
$ws.Cells.Item(11, 18).Value = [double]"705.9973187259881"
$ws.Cells.Item(11, 19).Value = [double]"0.1824260893354788"
$ws.Cells.Item(11, 20).Value = [double]"0.1824260893354788"

$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Vtn"
$ws.Cells.Item(12, 3).Value = "Itgb6"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = [double]"3"
$ws.Cells.Item(12, 6).Value = [double]"1"
$ws.Cells.Item(12, 7).Value = [double]"34.97741266666667"
$ws.Cells.Item(12, 8).Value = [double]"104.932238"
$ws.Cells.Item(12, 9).Value = [double]"0.5948251867999219"
$ws.Cells.Item(12, 10).Value = [double]"0.5948251867999219"
$ws.Cells.Item(12, 11).Value = [double]"3"
$ws.Cells.Item(12, 12).Value = [double]"1"
$ws.Cells.Item(12, 13).Value = [double]"4.789377333333333"
$ws.Cells.Item(12, 14).Value = [double]"14.368132"
$ws.Cells.Item(12, 15).Value = [double]"0.6549434327891901"
$ws.Cells.Item(12, 16).Value = [double]"0.6549434327891902"
$ws.Cells.Item(12, 17).Value = [double]"167.5200274043796"
$ws.Cells.Item(12, 18).Value = [double]"1507.680246639416"
$ws.Cells.Item(12, 19).Value = [double]"0.3895768497522121"
$ws.Cells.Item(12, 20).Value = [double]"0.3895768497522122"

$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Vtn"
$ws.Cells.Item(13, 3).Value = "Itgb6"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = [double]"3"
$ws.Cells.Item(13, 6).Value = [double]"1"
$ws.Cells.Item(13, 7).Value = [double]"34.97741266666667"
$ws.Cells.Item(13, 8).Value = [double]"104.932238"
$ws.Cells.Item(13, 9).Value = [double]"0.5948251867999219"
$ws.Cells.Item(13, 10).Value = [double]"0.5948251867999219"
$ws.Cells.Item(13, 11).Value = [double]"1"
$ws.Cells.Item(13, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(13, 13).Value = [double]"0.02281433333333334"
$ws.Cells.Item(13, 14).Value = [double]"0.068443"
$ws.Cells.Item(13, 15).Value = [double]"0.003119841422001868"
$ws.Cells.Item(13, 16).Value = [double]"0.003119841422001868"
$ws.Cells.Item(13, 17).Value = [double]"0.7979863517148891"
$ws.Cells.Item(13, 18).Value = [double]"7.181877165434002"
$ws.Cells.Item(13, 19).Value = [double]"0.001855760256628395"
$ws.Cells.Item(13, 20).Value = [double]"0.001855760256628395"

$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Vtn"
$ws.Cells.Item(14, 3).Value = "Itgb6"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = [double]"1"
$ws.Cells.Item(14, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(14, 7).Value = [double]"0.004706"
$ws.Cells.Item(14, 8).Value = [double]"0.014118"
$ws.Cells.Item(14, 9).Value = [double]"8.003014275976175E-05"
$ws.Cells.Item(14, 10).Value = [double]"8.003014275976175E-05"
$ws.Cells.Item(14, 11).Value = [double]"3"
$ws.Cells.Item(14, 12).Value = [double]"1"
$ws.Cells.Item(14, 13).Value = [double]"0.2577576666666667"
$ws.Cells.Item(14, 14).Value = [double]"0.7732730000000001"
$ws.Cells.Item(14, 15).Value = [double]"0.03524815007985697"
$ws.Cells.Item(14, 16).Value = [double]"0.03524815007985697"
$ws.Cells.Item(14, 17).Value = [double]"0.001213007579333334"
$ws.Cells.Item(14, 18).Value = [double]"0.010917068214"
$ws.Cells.Item(14, 19).Value = [double]"2.820914482908461E-06"
$ws.Cells.Item(14, 20).Value = [double]"2.820914482908461E-06"

$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Vtn"
$ws.Cells.Item(15, 3).Value = "Itgb6"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = [double]"1"
$ws.Cells.Item(15, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(15, 7).Value = [double]"0.004706"
$ws.Cells.Item(15, 8).Value = [double]"0.014118"
$ws.Cells.Item(15, 9).Value = [double]"8.003014275976175E-05"
$ws.Cells.Item(15, 10).Value = [double]"8.003014275976175E-05"
$ws.Cells.Item(15, 11).Value = [double]"3"
$ws.Cells.Item(15, 12).Value = [double]"1"
$ws.Cells.Item(15, 13).Value = [double]"2.242708666666667"
$ws.Cells.Item(15, 14).Value = [double]"6.728126"
$ws.Cells.Item(15, 15).Value = [double]"0.3066885757089511"
$ws.Cells.Item(15, 16).Value = [double]"0.3066885757089511"
$ws.Cells.Item(15, 17).Value = [double]"0.01055418698533333"
$ws.Cells.Item(15, 18).Value = [double]"0.094987682868"
$ws.Cells.Item(15, 19).Value = [double]"2.454433049677536E-05"
$ws.Cells.Item(15, 20).Value = [double]"2.454433049677536E-05"

$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Vtn"
$ws.Cells.Item(16, 3).Value = "Itgb6"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = [double]"1"
$ws.Cells.Item(16, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(16, 7).Value = [double]"0.004706"
$ws.Cells.Item(16, 8).Value = [double]"0.014118"
$ws.Cells.Item(16, 9).Value = [double]"8.003014275976175E-05"
$ws.Cells.Item(16, 10).Value = [double]"8.003014275976175E-05"
$ws.Cells.Item(16, 11).Value = [double]"3"
$ws.Cells.Item(16, 12).Value = [double]"1"
$ws.Cells.Item(16, 13).Value = [double]"4.789377333333333"
$ws.Cells.Item(16, 14).Value = [double]"14.368132"
$ws.Cells.Item(16, 15).Value = [double]"0.6549434327891901"
$ws.Cells.Item(16, 16).Value = [double]"0.6549434327891902"
$ws.Cells.Item(16, 17).Value = [double]"0.02253880973066667"
$ws.Cells.Item(16, 18).Value = [double]"0.202849287576"
$ws.Cells.Item(16, 19).Value = [double]"5.241521642568731E-05"
$ws.Cells.Item(16, 20).Value = [double]"5.241521642568732E-05"

$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "Vtn"
$ws.Cells.Item(17, 3).Value = "Itgb6"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = [double]"1"
$ws.Cells.Item(17, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(17, 7).Value = [double]"0.004706"
$ws.Cells.Item(17, 8).Value = [double]"0.014118"
$ws.Cells.Item(17, 9).Value = [double]"8.003014275976175E-05"
$ws.Cells.Item(17, 10).Value = [double]"8.003014275976175E-05"
$ws.Cells.Item(17, 11).Value = [double]"1"
$ws.Cells.Item(17, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(17, 13).Value = [double]"0.02281433333333334"
$ws.Cells.Item(17, 14).Value = [double]"0.068443"
$ws.Cells.Item(17, 15).Value = [double]"0.003119841422001868"
$ws.Cells.Item(17, 16).Value = [double]"0.003119841422001868"
$ws.Cells.Item(17, 17).Value = [double]"0.0001073642526666667"
$ws.Cells.Item(17, 18).Value = [double]"0.0009662782740000001"
$ws.Cells.Item(17, 19).Value = [double]"2.496813543906276E-07"
$ws.Cells.Item(17, 20).Value = [double]"2.496813543906277E-07"
